$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.487164258956909
$ws.Range("B1").Value = 4.344063758850098
$ws.Range("C1").Value = 3.230034112930298
$ws.Range("D1").Value = 0.8917667865753174
$ws.Range("E1").Value = 0.4708960652351379
